$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.116.89"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.668.83"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5213"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2637"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06227"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07497"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("D12").Value = "1.688.39"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.421"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("E14").Value = "  -4.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000007901"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.97%  "
$ws.Range("D17").Value = "26.178.86"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("E21").Value = "  -5.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.172"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1244"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.569"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06229"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.363"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("E30").Value = "  -4.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.479"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.423"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.626"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9950"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6039"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.702"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.123"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01603"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "1.074.46"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8654"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").Value = "1.816.37"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.946"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4251"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.941"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.17%  "
